# Append " – Done" (as two separate red-colored runs, matching the
# formatting already used elsewhere in this doc for resolved items) to
# the end of the "Observation date" max/min-date feedback bullet.
$d = $word.ActiveDocument

# wdColor value for hex C9211E (Word stores RGB as 0xBBGGRR)
$redColor = 1974729

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Can we block out future dates*") {
        $pEnd = $p.Range.End

        # Insert "Done" first, right before the paragraph mark, so it
        # inherits the paragraph's own run properties (incl. en-US lang).
        $r2 = $d.Range($pEnd - 1, $pEnd - 1)
        $r2.InsertAfter("Done")
        $r2.Font.Color = $redColor

        # Then insert "– " just before "Done", inheriting formatting from
        # the preceding "Can we block out future dates..." run.
        $r1 = $d.Range($pEnd - 1, $pEnd - 1)
        $r1.InsertAfter("– ")
        $r1.Font.Color = $redColor
        break
    }
}
